$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update underlying data values (diff: C5, E5, C6, E6, C7, E7, C8, C9, C10) ---
$ws.Range("C5").Value = 13297
$ws.Range("E5").Value = 1348
$ws.Range("C6").Value = 1409
$ws.Range("E6").Value = 635
$ws.Range("C7").Value = 1831
$ws.Range("E7").Value = 952
$ws.Range("C8").Value = 8301
$ws.Range("C9").Value = 808
$ws.Range("C10").Value = 1204

# --- Style changes ---
# J2: s=11 -> s=1 (drop the "applyFill" no-op fill flag; plain centered style)
# K5: s=11 -> s=1 (same target style as J2)
# Use H2 (already s=1) as the format donor.
$ws.Range("H2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("K5").PasteSpecial(-4122)

# K2: s=8 -> s=11 (drop the green highlight fill, back to plain/no-fill style)
# K8: s=8 -> s=11 (same)
# Use a donor cell that already carries the target "no fill" look.
$ws.Range("H5").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("K8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Selection / active cell (diff: E18 -> F19) ---
$ws.Range("F19").Select()

# --- Window position (diff: workbookView xWindow/yWindow) ---
$win = $excel.ActiveWindow
$win.Left = 5016
$win.Top = 2940

# Recalculate so dependent formulas (H5, H8, etc.) pick up new cached values
$excel.CalculateFull()
